# Fruta / hortaliza, semanal
# Insert a new weekly record at row 173 (pushing the existing rows 173-225
# down to 174-226) on the single data sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 173..225 down to 174..226, leaving a blank row 173.
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row 173 with the new weekly observation.
$ws.Cells.Item(173, 1).Value2  = 10
$ws.Cells.Item(173, 2).Value2  = 'Vega Modelo de Temuco'
$ws.Cells.Item(173, 3).Value2  = 'La Araucanía'
$ws.Cells.Item(173, 4).Value2  = 44809
$ws.Cells.Item(173, 5).Value2  = 9
$ws.Cells.Item(173, 6).Value2  = 'Fruta'
$ws.Cells.Item(173, 7).Value2  = 100104
$ws.Cells.Item(173, 8).Value2  = 'Frutos de pepita'
$ws.Cells.Item(173, 9).Value2  = 100104003
$ws.Cells.Item(173, 10).Value2 = 'Membrillo'
$ws.Cells.Item(173, 11).Value2 = 'Champion'
$ws.Cells.Item(173, 12).Value2 = 'Primera'
$ws.Cells.Item(173, 13).Value2 = 65
$ws.Cells.Item(173, 14).Value2 = 10000
$ws.Cells.Item(173, 15).Value2 = 10000
$ws.Cells.Item(173, 16).Value2 = 10000
$ws.Cells.Item(173, 17).Value2 = '$/bandeja 18 kilos granel'
$ws.Cells.Item(173, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(173, 19).Value2 = 556
$ws.Cells.Item(173, 20).Value2 = 18
